$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# "Remap sample from Field to FieldCollection" (#638)
#
# The getFirstOrNullObject/result=1 sample row that used to live under
# the "Field" class (row 32) is moved down to become a "FieldCollection"
# sample (row 35); the Field rows that followed it (code/parentBody/
# result) each shift up by one row to fill the gap. The final
# FieldCollection/items row (row 36) is unchanged.
# -----------------------------------------------------------------

# New content for rows 32-36, columns A-E (Class, Method/Prop/Rel Name,
# Member ID (methods only), SnippetIdIntheYAMLFile, MethodNameInTheSnippet)
$rows = @(
    @{ Row = 32; A = "Field";           B = "code";                  C = $null; D = "word-manage-fields"; E = "getFirstField" },
    @{ Row = 33; A = "Field";           B = "parentBody";            C = $null; D = "word-manage-fields"; E = "getParentBodyOfFirstField" },
    @{ Row = 34; A = "Field";           B = "result";                C = $null; D = "word-manage-fields"; E = "getFirstField" },
    @{ Row = 35; A = "FieldCollection"; B = "getFirstOrNullObject";  C = 1;     D = "word-manage-fields"; E = "getFirstField" },
    @{ Row = 36; A = "FieldCollection"; B = "items";                 C = $null; D = "word-manage-fields"; E = "getAllFields" }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Cells.Item($n, 1).Value = $r.A
    $ws.Cells.Item($n, 2).Value = $r.B
    if ($null -eq $r.C) {
        $ws.Cells.Item($n, 3).Value = ""
    } else {
        $ws.Cells.Item($n, 3).Value = $r.C
    }
    $ws.Cells.Item($n, 4).Value = $r.D
    $ws.Cells.Item($n, 5).Value = $r.E
}

# Normalize the row formatting (these rows previously used a redundant
# "applyNumberFormat" style variant; re-align them with the plain styles
# used by the rest of the table, e.g. row 31/37).
$ws.Range("A31").Copy() | Out-Null
$ws.Range("A32:A36").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("D31").Copy() | Out-Null
$ws.Range("B32:B36").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("D32:D36").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("E32:E36").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = $false

# Update the active selection / scroll position to match the edited area.
$ws.Range("E35").Select()
